$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Update the date field result: "22 de mayo de 2021" -> "24 de mayo de 2021"
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute('22 de mayo de 2021', $true, $false, $false, $false, $false, $true, 1, $false, '24 de mayo de 2021', 2)

# ------------------------------------------------------------------
# 2) "${model} ${brand}" -> "${model} " and move the _GoBack bookmark
#    to sit right after "${model} " (where "${brand}" used to be).
#    We do this by bookmarking the "${brand}" text first, then
#    deleting it -- the bookmark naturally collapses into place,
#    which also removes the old _GoBack bookmark (name is unique).
# ------------------------------------------------------------------
$modelPrefix = '${model} '
$rngModel = $d.Content
$null = $rngModel.Find.Execute('${model} ${brand}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$modelEnd = $rngModel.Start + $modelPrefix.Length

$rngBrand = $d.Range($modelEnd, $rngModel.End)
$null = $d.Bookmarks.Add("_GoBack", $rngBrand)

$rngBrand2 = $d.Range($modelEnd, $rngModel.End)
$rngBrand2.Text = ""

# ------------------------------------------------------------------
# 3) Split the lone "${Brand}" run into two runs: "${b" and "rand}"
#    Toggling Bold on/off after inserting text forces the engine to
#    split the run without leaving any residual formatting.
# ------------------------------------------------------------------
$rngBr = $d.Content
$null = $rngBr.Find.Execute('${Brand}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
$rngBr.Text = '${b'
$rngBr.Collapse(0)
$rngBr.InsertAfter('rand}')
$rngTail = $d.Range($rngBr.Start, $rngBr.Start + 5)
$rngTail.Bold = 1
$rngTail.Bold = 0
